# Auto-generated edit script for all_person_matches.xlsx
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "person matches"
$ws2 = $wb.Worksheets.Item(2)   # "grouped matches"

# --- Sheet 1 ("person matches"): append 13 new rows (91-103) ---
# Copy formatting from the last existing row (90) down into the new rows
# so the new rows pick up the same style (bold/bordered index column, etc.)
$ws1.Range("A90:H90").Copy()
$ws1.Range("A91:H103").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws1.Cells.Item(91, 1).Value = 42
$ws1.Cells.Item(91, 2).Value = 'Gö Chödrup'
$ws1.Cells.Item(91, 3).Value = '?'
$ws1.Cells.Item(91, 6).Value = 'eft:g-ch-drup'

$ws1.Cells.Item(92, 1).Value = 43
$ws1.Cells.Item(92, 2).Value = 'wang phab zhwun (wang phan zhun)'
$ws1.Cells.Item(92, 3).Value = '?'
$ws1.Cells.Item(92, 6).Value = 'eft:wang-phab-zhwun-wang-phan-zhun-'

$ws1.Cells.Item(93, 1).Value = 44
$ws1.Cells.Item(93, 2).Value = 'dge ba''i blo gros'
$ws1.Cells.Item(93, 3).Value = '?'
$ws1.Cells.Item(93, 6).Value = 'eft:dge-ba-i-blo-gros'

$ws1.Cells.Item(94, 1).Value = 45
$ws1.Cells.Item(94, 2).Value = 'rgya mtsho''i sde'
$ws1.Cells.Item(94, 3).Value = '?'
$ws1.Cells.Item(94, 6).Value = 'eft:rgya-mtsho-i-sde'

$ws1.Cells.Item(95, 1).Value = 46
$ws1.Cells.Item(95, 2).Value = 'Thönmi Sambhoṭa'
$ws1.Cells.Item(95, 3).Value = '?'
$ws1.Cells.Item(95, 6).Value = 'eft:th-nmi-sambhota'

$ws1.Cells.Item(96, 1).Value = 47
$ws1.Cells.Item(96, 2).Value = 'Tsultrim Gyaltsen'
$ws1.Cells.Item(96, 3).Value = '?'
$ws1.Cells.Item(96, 6).Value = 'eft:tsultrim-gyaltsen'

$ws1.Cells.Item(97, 1).Value = 48
$ws1.Cells.Item(97, 2).Value = 'Shang Buchikpa'
$ws1.Cells.Item(97, 3).Value = '?'
$ws1.Cells.Item(97, 6).Value = 'eft:shang-buchikpa'

$ws1.Cells.Item(98, 1).Value = 49
$ws1.Cells.Item(98, 2).Value = 'Sherap Ö'
$ws1.Cells.Item(98, 3).Value = '?'
$ws1.Cells.Item(98, 6).Value = 'eft:sherap-'

$ws1.Cells.Item(99, 1).Value = 50
$ws1.Cells.Item(99, 2).Value = 'Paṇḍita Dharmākara'
$ws1.Cells.Item(99, 3).Value = '?'
$ws1.Cells.Item(99, 6).Value = 'eft:pandita-dharmakara'

$ws1.Cells.Item(100, 1).Value = 51
$ws1.Cells.Item(100, 2).Value = 'Lotsāwa Zangkyong (bzang skyong)'
$ws1.Cells.Item(100, 3).Value = '?'
$ws1.Cells.Item(100, 6).Value = 'eft:lotsawa-zangkyong-bzang-skyong-'

$ws1.Cells.Item(101, 1).Value = 52
$ws1.Cells.Item(101, 2).Value = 'Nyen Lotsawa Darma Drak'
$ws1.Cells.Item(101, 3).Value = '?'
$ws1.Cells.Item(101, 6).Value = 'eft:nyen-lotsawa-darma-drak'

$ws1.Cells.Item(102, 1).Value = 53
$ws1.Cells.Item(102, 2).Value = 'Patsap Nyima Drak [?]'
$ws1.Cells.Item(102, 3).Value = '?'
$ws1.Cells.Item(102, 6).Value = 'eft:patsap-nyima-drak-'

$ws1.Cells.Item(103, 1).Value = 54
$ws1.Cells.Item(103, 2).Value = 'vajrvisramitra'
$ws1.Cells.Item(103, 3).Value = '?'
$ws1.Cells.Item(103, 6).Value = 'eft:vajrvisramitra'

# --- Sheet 2 ("grouped matches"): refresh the BDRC-ID / 84000-ID groupings ---
$ws2.Cells.Item(2, 2).Value = 'P8213'
$ws2.Cells.Item(2, 3).Value = '{''eft:t-vidyakarasimha'', ''eft:vidyakarasimha''}'
$ws2.Cells.Item(3, 2).Value = 'P8268'
$ws2.Cells.Item(3, 3).Value = '{''eft:buddhaprabha''}'
$ws2.Cells.Item(4, 2).Value = 'P0TMP080'
$ws2.Cells.Item(4, 3).Value = '{''eft:hwa-shang-zab-mo''}'
$ws2.Cells.Item(5, 2).Value = 'P4242'
$ws2.Cells.Item(5, 3).Value = '{''eft:sherab-lekpa''}'
$ws2.Cells.Item(6, 2).Value = 'P8273'
$ws2.Cells.Item(6, 3).Value = '{''eft:rinchen-tso'', ''eft:rin-chen-tsho''}'
$ws2.Cells.Item(7, 2).Value = 'P3456'
$ws2.Cells.Item(7, 3).Value = '{''eft:tshul-khrims-rgyal-ba''}'
$ws2.Cells.Item(8, 2).Value = 'P8267'
$ws2.Cells.Item(8, 3).Value = '{''eft:vijayasila''}'
$ws2.Cells.Item(9, 2).Value = 'P8260'
$ws2.Cells.Item(9, 3).Value = '{''eft:dpal-dbyangs''}'
$ws2.Cells.Item(10, 2).Value = 'P8266'
$ws2.Cells.Item(10, 3).Value = '{''eft:dharmatasila'', ''eft:ch-nyi-tsultrim''}'
$ws2.Cells.Item(11, 2).Value = 'P0TMP104'
$ws2.Cells.Item(11, 3).Value = '{''eft:punyasambhava''}'
$ws2.Cells.Item(12, 2).Value = 'P8217'
$ws2.Cells.Item(12, 3).Value = '{''eft:jnanagarbha'', ''eft:t-jnanagarbha''}'
$ws2.Cells.Item(13, 2).Value = 'P4258'
$ws2.Cells.Item(13, 3).Value = '{''eft:dpal-byor''}'
$ws2.Cells.Item(14, 2).Value = 'https://lod.dila.edu.tw/resource.php?id=A000089'
$ws2.Cells.Item(14, 3).Value = '{''eft:siladharma''}'
$ws2.Cells.Item(15, 2).Value = '?'
$ws2.Cells.Item(15, 3).Value = '{''Shang Buchikpa'', ''Patsap Nyima Drak [?]'', ''vajrvisramitra'', ''Thönmi Sambhoṭa'', ''Sherap Ö'', ''wang phab zhwun (wang phan zhun)'', ''Gö Chödrup'', "dge ba''i blo gros", ''Lotsāwa Zangkyong (bzang skyong)'', ''Tsultrim Gyaltsen'', ''Paṇḍita Dharmākara'', ''Nyen Lotsawa Darma Drak'', ''eft:sakyasena'', "rgya mtsho''i sde"}'
$ws2.Cells.Item(16, 2).Value = 'P2956'
$ws2.Cells.Item(16, 3).Value = '{''eft:krsnapandita''}'
$ws2.Cells.Item(17, 2).Value = 'P0RK8'
$ws2.Cells.Item(17, 3).Value = '{''eft:dharmapala''}'
$ws2.Cells.Item(18, 2).Value = 'P3709'
$ws2.Cells.Item(18, 3).Value = '{''eft:phakpa-sherab''}'
$ws2.Cells.Item(19, 2).Value = 'P0TMPT007'
$ws2.Cells.Item(19, 3).Value = '{''eft:rnam-par-mi-rtog-pa''}'
$ws2.Cells.Item(20, 2).Value = 'P4CZ16780'
$ws2.Cells.Item(20, 3).Value = '{''eft:manjusrigarbha''}'
$ws2.Cells.Item(21, 2).Value = 'P3379'
$ws2.Cells.Item(21, 3).Value = '{''eft:dipamkarasrijnana'', ''eft:dipamkara-srijnana''}'
$ws2.Cells.Item(22, 2).Value = 'P8211'
$ws2.Cells.Item(22, 3).Value = '{''eft:vidyakaraprabha''}'
$ws2.Cells.Item(23, 2).Value = 'P8220'
$ws2.Cells.Item(23, 3).Value = '{''eft:devacandra''}'
$ws2.Cells.Item(24, 2).Value = 'P3214'
$ws2.Cells.Item(24, 3).Value = '{''eft:danasila''}'
$ws2.Cells.Item(25, 2).Value = 'P753'
$ws2.Cells.Item(25, 3).Value = '{''eft:rin-chen-bzang-po''}'
$ws2.Cells.Item(26, 2).Value = 'P8171'
$ws2.Cells.Item(26, 3).Value = '{''eft:dharmasribhadra''}'
$ws2.Cells.Item(27, 2).Value = 'P3285'
$ws2.Cells.Item(27, 3).Value = '{''eft:sakya-yesh-''}'
$ws2.Cells.Item(28, 2).Value = 'P2548'
$ws2.Cells.Item(28, 3).Value = '{''eft:prajnavarman'', ''eft:prajnavarma''}'
$ws2.Cells.Item(29, 2).Value = 'P8228'
$ws2.Cells.Item(29, 3).Value = '{''eft:surendrabodhi''}'
$ws2.Cells.Item(30, 2).Value = 'P8263'
$ws2.Cells.Item(30, 3).Value = '{''eft:leki-d-''}'
$ws2.Cells.Item(31, 2).Value = 'P5651'
$ws2.Cells.Item(31, 3).Value = '{''eft:pa-tshab-nyi-ma-grags''}'
$ws2.Cells.Item(32, 2).Value = 'P0TMP092'
$ws2.Cells.Item(32, 3).Value = '{''eft:anandasri-s-''}'
$ws2.Cells.Item(33, 2).Value = 'P4259'
$ws2.Cells.Item(33, 3).Value = '{''eft:palgyi-lh-npo'', ''eft:dpal-gyi-lhun-po'', ''eft:ban-de-dpal-gyi-lhun-po''}'
$ws2.Cells.Item(34, 2).Value = 'P4CZ16819'
$ws2.Cells.Item(34, 3).Value = '{''eft:sakyaprabha''}'
$ws2.Cells.Item(35, 2).Value = 'P1KG8854'
$ws2.Cells.Item(35, 3).Value = '{''eft:surendrabodhi'', ''eft:silendrabodhi'', ''eft:srilendrabodhi''}'
$ws2.Cells.Item(36, 2).Value = 'P8222'
$ws2.Cells.Item(36, 3).Value = '{''eft:jnanasiddhi'', ''eft:jnanasidhi''}'
$ws2.Cells.Item(37, 2).Value = 'P8245'
$ws2.Cells.Item(37, 3).Value = '{''eft:buddhakaravarma''}'
$ws2.Cells.Item(38, 2).Value = 'P2637'
$ws2.Cells.Item(38, 3).Value = '{''eft:trakpa-gyaltsen''}'
$ws2.Cells.Item(39, 2).Value = 'P8205'
$ws2.Cells.Item(39, 3).Value = '{''eft:zhang-yesh-d-'', ''eft:ye-shes-sde'', ''eft:yesh-d-'', ''eft:band-yesh-d-'', ''eft:band-yesh-de'', ''eft:yesh-d-ye-shes-sde-''}'
$ws2.Cells.Item(40, 2).Value = 'P8249'
$ws2.Cells.Item(40, 3).Value = '{''eft:dharmakara''}'
$ws2.Cells.Item(41, 2).Value = 'P8093'
$ws2.Cells.Item(41, 3).Value = '{''eft:kamalagupta''}'
$ws2.Cells.Item(42, 2).Value = 'P8206'
$ws2.Cells.Item(42, 3).Value = '{''eft:celu''}'
$ws2.Cells.Item(43, 2).Value = 'P8183'
$ws2.Cells.Item(43, 3).Value = '{''eft:klu-i-rgyal-mtshan'', ''eft:cog-ro-klu-i-rgyal-mtshan''}'
$ws2.Cells.Item(44, 2).Value = 'P8261'
$ws2.Cells.Item(44, 3).Value = '{''eft:munivarma'', ''eft:munivarman''}'
$ws2.Cells.Item(45, 2).Value = 'P8209'
$ws2.Cells.Item(45, 3).Value = '{''eft:jinamitra'', ''eft:jinamitra-k-'', ''eft:dzi-na-mi-tra-k-''}'
$ws2.Cells.Item(46, 2).Value = 'P8182'
$ws2.Cells.Item(46, 3).Value = '{''eft:kawa-paltsek-under-the-name-paltsek-raksita-'', ''eft:dpal-brtsegs'', ''eft:ska-ba-dpal-brtsegs'', ''eft:paltsek'', ''eft:ban-de-dpal-brtsegs''}'
$ws2.Cells.Item(47, 2).Value = 'P4263'
$ws2.Cells.Item(47, 3).Value = '{''eft:dge-ba-dpal''}'
$ws2.Cells.Item(48, 2).Value = 'P8151'
$ws2.Cells.Item(48, 3).Value = '{''eft:gayadhara''}'
$ws2.Cells.Item(49, 2).Value = 'P0TMP098'
$ws2.Cells.Item(49, 3).Value = '{''eft:jinavara''}'
$ws2.Cells.Item(50, 2).Value = 'P8219'
$ws2.Cells.Item(50, 3).Value = '{''eft:visuddhasimha''}'
$ws2.Cells.Item(51, 2).Value = 'P8265'
$ws2.Cells.Item(51, 3).Value = '{''eft:ratnaraksita''}'
$ws2.Cells.Item(52, 2).Value = 'P4255'
$ws2.Cells.Item(52, 3).Value = '{''eft:ye-shes-snying-po'', ''eft:t-jnanagarbha'', ''eft:yesh-nyingpo''}'
$ws2.Cells.Item(53, 2).Value = 'P4CZ15137'
$ws2.Cells.Item(53, 3).Value = '{''eft:kumarakalasa''}'
$ws2.Cells.Item(54, 2).Value = 'P00KG07267'
$ws2.Cells.Item(54, 3).Value = '{''eft:sarvajnadeva'', ''eft:sarvanyadeva''}'
$ws2.Cells.Item(55, 2).Value = 'P8269'
$ws2.Cells.Item(55, 3).Value = '{''eft:dgon-gling-rma''}'

Write-Output "edit applied"
